$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Smoke Particle system row (row 12): mark the item as DONE and bump its score
$ws.Range("D12").Value = 10
$ws.Range("E12").Value = "DONE"
# Match the green "DONE" font styling used by the other completed rows (e.g. E4)
$ws.Range("E12").Font.Color = $ws.Range("E4").Font.Color
$ws.Range("F12").Value = "better blend mode?"

# Restore the active selection on the frozen pane to F17
$ws.Range("F17").Select()

$wb.Save()
